# feat: add 2022-Q4 data
#
# 1) Insert a brand-new "2022-Q4" sheet (fund-holdings detail) right after
#    the "总计" (totals) sheet, pushing all the quarter sheets one slot later.
# 2) Update the "总计" sheet: insert a new summary row for 2022-Q4 right
#    after the header, and append a duplicate-valued row for 2021-Q2 at
#    the bottom (mirrors the upstream data export's index bookkeeping).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q4" worksheet after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Header row (bold/centered style matches the other quarter sheets, so
# copy formatting from an existing quarter sheet's header row first).
$templateSheet = $wb.Worksheets.Item(3)
$templateSheet.Range("A1:H2").Copy($q4Sheet.Range("A1"))

$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$q4Sheet.Range("A2").Value = 0

# Text-like numeric columns must stay text (leading zeros / exact string
# form matter), so force text format before assigning.
$q4Sheet.Range("B2:G2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "010404"
$q4Sheet.Range("C2").Value = "博道盛利6个月持有期混合"
$q4Sheet.Range("D2").Value = "1.07"
$q4Sheet.Range("E2").Value = "41.08"
$q4Sheet.Range("F2").Value = "0.36"
$q4Sheet.Range("G2").Value = "0.0039"
$q4Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Insert a fresh row 2 for the 2022-Q4 summary, pushing existing quarters
# down by one row.
$ws.Rows.Item(2).Insert()
$ws.Range("A3").Copy($ws.Range("A2"))
$ws.Range("B2:D2").ClearFormats()
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

# Renumber the running index in column A for the shifted rows.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# Append a new bottom row (2021-Q2), duplicating the last row's values.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("B6:D6").Copy($ws.Range("B7"))
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "2021-Q2"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.01
